$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 86 with inventory item data
$ws.Range("A86").Value = "L8SEY6"
$ws.Range("B86").Value = "Almohadilla+Chip Epson T671"
$ws.Range("C86").Value = "WF 6090 6091 6093 6590 6593 8010 8090 8093 8510 8590 8591 8593"
$ws.Range("D86").Value = 35000
$ws.Range("E86").Value = 200000
$ws.Range("F86").Value = 10
$ws.Range("G86").Value = 0
$ws.Range("H86").Formula = "=(E86-D86)*G86"
$ws.Range("I86").Formula = "=D86*F86"
$ws.Range("J86").Value = 350000

# Match number formatting style of the rest of the price/ganancia/inversion columns
$ws.Range("D86").NumberFormat = $ws.Range("D85").NumberFormat
$ws.Range("E86").NumberFormat = $ws.Range("E85").NumberFormat
$ws.Range("H86").NumberFormat = $ws.Range("H85").NumberFormat
$ws.Range("I86").NumberFormat = $ws.Range("I85").NumberFormat
